$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format specific cells as Text so numeric-looking strings with
# trailing zeros (e.g. "551.50") are preserved exactly as text, matching
# the original inline-string cell contents.
foreach ($addr in @("D5", "D19", "D35", "D39", "D45")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '59.896.70'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '2.421.33'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '551.50'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").Value = '137.46'
$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("E8").Value = '  +2.52%  '

$ws.Range("E9").Value = '  -2.19%  '

$ws.Range("D10").Value = '5.68'
$ws.Range("E10").Value = '  -2.75%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.147'
$ws.Range("E11").Value = '  -2.26%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '0.354'
$ws.Range("E12").Value = '  -1.86%  '

$ws.Range("D13").Value = '25.48'
$ws.Range("E13").Value = '  +3.84%  '

$ws.Range("D14").Value = '2.853.98'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '59.826.78'
$ws.Range("E15").Value = '  +1.00%  '

$ws.Range("E16").Value = '  -1.52%  '

$ws.Range("D17").Value = '2.427.22'
$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").Value = '4.40'
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").Value = '330.78'
$ws.Range("E20").Value = '  -1.49%  '

$ws.Range("E21").Value = '  -3.99%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").Value = '66.38'
$ws.Range("E23").Value = '  +2.87%  '

$ws.Range("E24").Value = '  +1.15%  '

$ws.Range("D25").Value = '8.75'
$ws.Range("E25").Value = '  +3.76%  '

$ws.Range("E26").Value = '  +0.41%  '

$ws.Range("E27").Value = '  +1.89%  '

$ws.Range("E28").Value = '  +2.21%  '

$ws.Range("E29").Value = '  -0.83%  '

$ws.Range("D30").Value = '168.66'
$ws.Range("E30").Value = '  -1.45%  '

$ws.Range("D31").Value = '6.13'
$ws.Range("E31").Value = '  -1.80%  '

$ws.Range("D32").Value = '18.66'
$ws.Range("E32").Value = '  -0.41%  '

$ws.Range("E33").Value = '  +0.58%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = '1.30'
$ws.Range("E35").Value = '  +1.90%  '

$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("D39").Value = '39.60'
$ws.Range("E39").Value = '  -1.89%  '

$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("D41").Value = '313.84'
$ws.Range("E41").Value = '  +6.21%  '

$ws.Range("D42").Value = '3.67'
$ws.Range("E42").Value = '  -2.00%  '

$ws.Range("D43").Value = '139.34'
$ws.Range("E43").Value = '  -1.78%  '

$ws.Range("D44").Value = '0.0967'
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").Value = '0.0520'
$ws.Range("E45").Value = '  -0.60%  '

$ws.Range("D46").Value = '19.51'
$ws.Range("E46").Value = '  +2.87%  '

$ws.Range("D47").Value = '0.577'
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0225'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("B49").Value = 'Polygon'
$ws.Range("C49").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D49").Value = '0.395'
$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("E50").Value = '  -1.08%  '

$ws.Range("D51").Value = '11.05'
$ws.Range("E51").Value = '  +0.13%  '
